$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35 updates ---
$ws.Range("B35").Value = "[57, 343]"
$ws.Range("C35").Value = 0.3575
$ws.Range("D35").Value = "[1, 190, 1, 208]"
$ws.Range("E35").Value = 0.50625
$ws.Range("F35").Value = "[2, 6, 4, 12, 45, 9, 7, 1, 26, 13, 1, 2, 5, 3, 16, 8, 4, 7, 7, 9, 7, 4, 12, 8, 9, 4, 3, 7, 1, 3, 1, 7, 2, 5, 19, 13, 5, 8, 17, 1, 8, 1, 9, 34, 3, 4, 1, 8, 7, 1, 1]"
$ws.Range("G35").Value = 0.48
$ws.Range("H35").Value = "[1, 90, 1, 182, 118, 2, 4, 1, 1]"
$ws.Range("I35").Value = 0.697778

# --- Row 37 updates ---
$ws.Range("B37").Value = "[4, 0]"
$ws.Range("D37").Value = "[0, 0, 0, 4]"
$ws.Range("F37").Value = "[0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 3, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0]"
$ws.Range("G37").Value = 0.970588
$ws.Range("H37").Value = "[0, 0, 0, 0, 4, 0, 0, 0, 0]"

# --- Row 38 updates ---
$ws.Range("D38").Value = "[0, 3, 0, 0]"
$ws.Range("F38").Value = "[0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 3, 0, 0, 0, 0, 0]"

# --- Row 39 updates ---
$ws.Range("F39").Value = "[0, 0, 0, 0, 0, 3, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0]"

# --- Row 44 updates ---
$ws.Range("D44").Value = "[0, 2, 0, 0]"
$ws.Range("F44").Value = "[0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0]"
$ws.Range("H44").Value = "[0, 0, 0, 2, 0, 0, 0, 0, 0]"

# --- Row 52 updates ---
$ws.Range("D52").Value = "[0, 0, 0, 1]"

# --- New row 53, inserted after row 52 (a new community entry) ---
# Copy formatting of row 52 down to the new row 53 first so styles (fills on C/E/G/I) match.
$ws.Range("A52:I52").Copy() | Out-Null
$ws.Range("A53:I53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A53").Value = "Asynchronous label propagation  community 18"
$ws.Range("B53").Value = "[1, 0]"
$ws.Range("C53").Value = 0.5
$ws.Range("D53").Value = "[0, 1, 0, 0]"
$ws.Range("E53").Value = 0.75
$ws.Range("F53").Value = "[0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0]"
$ws.Range("G53").Value = 0.980392
$ws.Range("H53").Value = "[0, 0, 0, 0, 1, 0, 0, 0, 0]"
$ws.Range("I53").Value = 0.888889
